$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (scheduled GitHub Actions run): update the Price (D)
# and Volume(1h) (E) columns for every coin row with freshly fetched figures.
#
# D-column prices are plain text in the sheet (e.g. "29.380.59", trailing-zero
# values like "0.7147", tiny decimals like "0.000007900"), so each Price write
# is forced to text via NumberFormat "@" before assignment (otherwise Excel
# autoconverts numeric-looking text to a Number and mangles/ truncates it), then
# the style is reset back to "Normal" so no lingering text-format style remains
# applied to the cell (matches the original unstyled inline-string cells).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.380.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.874.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7147"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07826"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3082"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.89%  "
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.856.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7236"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.267"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.418.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.892"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "243.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007900"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.121.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.935"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.89%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1552"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.987"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("E29").Value = "  -4.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.481"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.368"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.110"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05265"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.921"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.196"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7190"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.678"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01859"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.212.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.712"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9063"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.067"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5342"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("E47").Value = "  +4.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.757"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.909"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4305"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.214"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.40%  "
